# EnemyDB: add new enemy-species rows (Spider, Skeleton, Orc, Golem, Bat,
# Dragon, MonsterPlant) to the Entities sheet, and set the AgentType-ish
# rank/size columns for each, per commit message:
#   "Enemy DB에 다른 종족들 ID별로 추가" / "Enemy별로 AgentType 설정"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# id, name, rank, size, E, attackRange, meleeAtk, magicAtk, def, attackSpeed, moveSpeed
$newRows = @(
    @(30200000, "Spider",       "Normal", "Small",  20, 1.6, 3, 0, 3, 0.5, 2),
    @(30300000, "Skeleton",     "Normal", "Small",  20, 2.3, 3, 0, 3, 0.5, 2),
    @(30400001, "Orc",          "Normal", "Medium", 20, 3,   3, 0, 3, 0.5, 2),
    @(30700001, "Golem",        "Normal", "Medium", 20, 2.5, 3, 0, 3, 0.5, 2),
    @(30600000, "Bat",          "Normal", "Small",  20, 2,   3, 0, 3, 0.5, 2),
    @(30700000, "Dragon",       "Normal", "Small",  20, 3,   3, 0, 3, 0.5, 2),
    @(30800000, "MonsterPlant", "Normal", "Small",  20, 2,   3, 0, 3, 0.5, 2)
)

$r = 8
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
    $ws.Cells.Item($r, 8).Value2 = $row[7]
    $ws.Cells.Item($r, 9).Value2 = $row[8]
    $ws.Cells.Item($r, 10).Value2 = $row[9]
    $ws.Cells.Item($r, 11).Value2 = $row[10]
    $r = $r + 1
}

# Column B now holds longer species names ("MonsterPlant") - widen/refit it,
# splitting it away from column A's width (which stays 10.5 / bestFit).
$ws.Columns.Item(2).ColumnWidth = 12.43

# Reflect the last-edited cell as the active selection, like the author left it.
[void]$ws.Range("F12").Select()
